$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '29.321.12'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.872.82'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.76'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3106'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07766'
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.07'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08397'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = '1.873.72'
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.235'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7108'
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.07'
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '29.332.44'
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.059'
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008188'
$ws.Range("E18").Value = '  +4.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.47'
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '2.120.24'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.750'
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1597'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.69'
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.023'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.48'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.402'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.318'
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.286'
$ws.Range("E32").Value = '  -3.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05291'
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.936'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.175'
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7451'
$ws.Range("E36").Value = '  -6.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").Value = '1.218.54'
$ws.Range("E39").Value = '  +4.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.723'
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.512'
$ws.Range("E41").Value = '  +4.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.79'
$ws.Range("E42").Value = '  +7.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8863'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.44'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '2.018.91'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.798'
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4313'
$ws.Range("E51").Value = '  +0.96%  '
